$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Journal Dev" tracker with the latest status -------------
# Row 6 = "Faire ReadMe.md" (date 44992) -> now completed: mark "Fait"
# (green fill, same as the other completed rows) and add the remark.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Value = "ReadMe.md réalisé dans le temps que je me suis donné"

# Row 7 = "Faire Maquette" -> now completed: mark "Fait" (same green fill)
# while keeping its original (thin, all-around) border, and add the remark.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").Borders(7).LineStyle = 1
$ws.Range("C7").Borders(7).Weight = 2
$ws.Range("C7").Borders(8).LineStyle = 1
$ws.Range("C7").Borders(8).Weight = 2
$ws.Range("C7").Borders(9).LineStyle = 1
$ws.Range("C7").Borders(9).Weight = 2
$ws.Range("C7").Borders(10).LineStyle = 1
$ws.Range("C7").Borders(10).Weight = 2
$ws.Range("E7").Value = "Maquette page requêtes faites "

$excel.CutCopyMode = 0

# --- Drop the unused, style-only F:G columns (no data ever lived there) --
$ws.Range("F1:G11").Clear()

# --- Leave the cursor where the author left it after the update ----------
$ws.Range("E7").Select() | Out-Null
